$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 144 (id=142, Northeast United vs ? on 2024-02-21): fill in the final
# score / result (FTHG, FTAG, FTR) that was not yet known, and refresh the
# closing odds columns (N:AC) now that the match has finished.
# ---------------------------------------------------------------------------
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = "A"

$ws.Range("N144").Value = 1.6
$ws.Range("O144").Value = 4
$ws.Range("P144").Value = 5.5
$ws.Range("Q144").Value = -1
$ws.Range("R144").Value = 2.025
$ws.Range("S144").Value = 1.825
$ws.Range("U144").Value = 1.825
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 4.5
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.825
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 1.025

# ---------------------------------------------------------------------------
# Row 145 (id=143, 2024-02-22): refresh the closing odds columns (N:S) with
# the latest market prices.
# ---------------------------------------------------------------------------
$ws.Range("N145").Value = 2.05
$ws.Range("O145").Value = 3.4
$ws.Range("P145").Value = 3.5
$ws.Range("Q145").Value = -0.25
$ws.Range("R145").Value = 1.775
$ws.Range("S145").Value = 2.025

# ---------------------------------------------------------------------------
# Row 146 (new match, id=144, 2024-02-23): Chennaiyin FC vs Mumbai City FC.
# Copy formatting from the row above for the styled columns (A = id,
# E = Date), then populate every field with the new fixture's data.
# ---------------------------------------------------------------------------
$ws.Range("A145").Copy($ws.Range("A146"))
$ws.Range("E145").Copy($ws.Range("E146"))

$ws.Range("A146").Value = 144
$ws.Range("B146").Value = 7749870
$ws.Range("C146").Value = "India Super League"
$ws.Range("D146").Value = "India Super League"
$ws.Range("E146").Value = 45345.45833333334
$ws.Range("F146").Value = "Chennaiyin FC"
$ws.Range("G146").Value = "Mumbai City FC"

$ws.Range("K146").Value = 3.6
$ws.Range("L146").Value = 3.75
$ws.Range("M146").Value = 1.909
$ws.Range("N146").Value = 3.8
$ws.Range("O146").Value = 3.75
$ws.Range("P146").Value = 1.85
$ws.Range("Q146").Value = 0.5
$ws.Range("R146").Value = 1.95
$ws.Range("S146").Value = 1.85
$ws.Range("T146").Value = 2.75
$ws.Range("U146").Value = 1.925
$ws.Range("V146").Value = 1.875
$ws.Range("W146").Value = 0
$ws.Range("X146").Value = 0
$ws.Range("Y146").Value = 0
$ws.Range("Z146").Value = 0
$ws.Range("AA146").Value = 0
